# Updates to execute RAD Extension Payment Type.
#
# Simulates the Katalon RAD test runner having executed every PaymentType
# except "Extension Payments" in both sheets:
#   - For every row whose "Execute" flag (column C) was consumed by the run,
#     clear the Execute cell (it is fully removed from the row, matching
#     Apache POI's "cell no longer written" behaviour).
#   - The "Extension Payments" rows keep their Execute flag (column C) since
#     they are queued up to run next.
#   - Rows that already finished running their Extension Payments case get a
#     fresh Result/Date (sheet1: rows 11-12 flip Fail -> Pass with new
#     timestamps; sheet2: rows 10-13 stay Pass but get new timestamps).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: FEINmismatch
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("FEINmismatch")

# Update the two "Extension Payments" rows: Result Fail -> Pass, fresh Date.
$ws1.Range("A11").Value = "Pass"
$ws1.Range("B11").Value = "Wed Mar 20 23:01:00 EDT 2024"
$ws1.Range("A12").Value = "Pass"
$ws1.Range("B12").Value = "Wed Mar 20 23:01:13 EDT 2024"

# Remove the "Execute" (column C) flag from every row that already ran,
# leaving it only on the still-queued "Extension Payments" rows (11-12).
$ws1.Range("C2:C10").Clear()
$ws1.Range("C13:C30").Clear()

# Move the selection to the next batch of rows queued for execution.
[void]$ws1.Range("C13:C30").Select()

# ---------------------------------------------------------------------
# Sheet 2: FEINSSNmismatch
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("FEINSSNmismatch")

# The four "Extension Payments" rows were already Pass; just refresh Date.
$ws2.Range("B10").Value = "Wed Mar 20 23:01:29 EDT 2024"
$ws2.Range("B11").Value = "Wed Mar 20 23:01:40 EDT 2024"
$ws2.Range("B12").Value = "Wed Mar 20 23:01:52 EDT 2024"
$ws2.Range("B13").Value = "Wed Mar 20 23:02:03 EDT 2024"

# Remove the "Execute" (column C) flag from every row that already ran,
# leaving it only on the "Extension Payments" rows (10-13).
$ws2.Range("C2:C9").Clear()
$ws2.Range("C14:C19").Clear()

# Move the selection to the next batch of rows queued for execution, and
# make this sheet the active tab (matches the saved workbook view).
[void]$ws2.Activate()
[void]$ws2.Range("C14:C19").Select()
